# Add the magic_item_table sheet at the end of the workbook (after the last existing sheet)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "magic_item_table"

# Row 1: headers
$newSheet.Cells.Item(1, 1).Value = "n_id"
$newSheet.Cells.Item(1, 2).Value = "u_uid_id"
$newSheet.Cells.Item(1, 3).Value = "s_title"
$newSheet.Cells.Item(1, 4).Value = "s_description"
$newSheet.Cells.Item(1, 5).Value = "n_category_id"
$newSheet.Cells.Item(1, 6).Value = "n_rarity_id"
$newSheet.Cells.Item(1, 7).Value = "n_price"
$newSheet.Cells.Item(1, 8).Value = "b_attunement"
$newSheet.Cells.Item(1, 9).Value = "d_last_update"
$newSheet.Cells.Item(1, 10).Value = "t_write"

# Row 2: types
$newSheet.Cells.Item(2, 1).Value = "integer"
$newSheet.Cells.Item(2, 2).Value = "uuid"
$newSheet.Cells.Item(2, 3).Value = "string"
$newSheet.Cells.Item(2, 4).Value = "string"
$newSheet.Cells.Item(2, 5).Value = "integer"
$newSheet.Cells.Item(2, 6).Value = "integer"
$newSheet.Cells.Item(2, 7).Value = "number .2"
$newSheet.Cells.Item(2, 8).Value = "boolean"
$newSheet.Cells.Item(2, 9).Value = "timestamp"
$newSheet.Cells.Item(2, 10).Value = "timestamp"

# Make the new sheet the active sheet/tab
$newSheet.Activate()
